$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Note messages: drop "successfully" wording (rows 50-51 keep their row,
#     only the message text changes) ---
$ws.Cells.Item(50, 3).Value = "The note has been updated."
$ws.Cells.Item(51, 3).Value = "The note has been created."

# --- Row 53 used to be "901 Authentification Not logged in." -> becomes the
#     relocated "607 Note The note has been deleted." row ---
$ws.Cells.Item(53, 1).Value = 607
$ws.Cells.Item(53, 2).Value = "Note"
$ws.Cells.Item(53, 3).Value = "The note has been deleted."
$ws.Cells.Item(53, 4).Value = "Response"

# --- Authentification block renumbered 900-903 and re-worded ---
$ws.Cells.Item(54, 1).Value = 900
$ws.Cells.Item(54, 2).Value = "Authentification"
$ws.Cells.Item(54, 3).Value = "Not logged in."
$ws.Cells.Item(54, 4).Value = "Response"

$ws.Cells.Item(55, 1).Value = 901
$ws.Cells.Item(55, 2).Value = "Authentification"
$ws.Cells.Item(55, 3).Value = "Access is not authorized."
$ws.Cells.Item(55, 4).Value = "Response"

$ws.Cells.Item(56, 1).Value = 902
$ws.Cells.Item(56, 2).Value = "Authentification"
$ws.Cells.Item(56, 3).Value = "Failed to log in."
$ws.Cells.Item(56, 4).Value = "Response"

$ws.Cells.Item(57, 1).Value = 903
$ws.Cells.Item(57, 2).Value = "Authentification"
$ws.Cells.Item(57, 3).Value = "Account is already exist."
$ws.Cells.Item(57, 4).Value = "Response"

# --- New Authentification rows for parent registration / admin handling ---
$ws.Cells.Item(58, 1).Value = 904
$ws.Cells.Item(58, 2).Value = "Authentification"
$ws.Cells.Item(58, 3).Value = "Not a valid e-mail address."
$ws.Cells.Item(58, 4).Value = "Response"

$ws.Cells.Item(59, 1).Value = 905
$ws.Cells.Item(59, 2).Value = "Authentification"
$ws.Cells.Item(59, 3).Value = "New parent has been registered."
$ws.Cells.Item(59, 4).Value = "Response"

$ws.Cells.Item(60, 1).Value = 906
$ws.Cells.Item(60, 2).Value = "Authentification"
$ws.Cells.Item(60, 3).Value = "Not an admin account."
$ws.Cells.Item(60, 4).Value = "Response"

# --- MySQL 42S22 row moves to the very end (row 61) ---
$ws.Cells.Item(61, 1).Value = "42S22"
$ws.Cells.Item(61, 2).Value = "MySQL"
$ws.Cells.Item(61, 3).Value = "Unknown column in a database table."
$ws.Cells.Item(61, 4).Value = "Response"

# --- View: selection moves to C60, scrolled so row 43 is at the top ---
$ws.Range("C60").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
